# "nearly done with movement"
# - fill in the (previously blank) D2:D4 cells with an explicit 0, matching
#   the rest of the D column
# - move the active selection to L7 (was J14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0

$ws.Range("L7").Select() | Out-Null
